$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "ChartSheet_0"
$ws.Range("A11").Value = "gia tri moi "
$ws.Range("O19").Value = ""
